$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New strings must be entered in this exact order so the shared-string
# table ends up matching the original authoring session.

# Row 42 - headers for the new "Pair-integration" table
$ws.Range("A42").Value = "Pair-integration"
$ws.Range("B42").Value = "Duration"
$ws.Range("C42").Value = "Errors introduced into code"

# Row 43
$ws.Range("A43").Value = "Only one pair - everything right"
$ws.Range("B43").Value = 43
$ws.Range("C43").Value = 4

# Row 44
$ws.Range("A44").Value = "Two pairs - everything right"

# Row 45
$ws.Range("A45").Value = "All pairs - everything right"
$ws.Range("B45").Value = 17
$ws.Range("C45").Value = 3

# Row 47 (filled in before row 46)
$ws.Range("A47").Value = "No coding standard"

# Row 48
$ws.Range("A48").Value = "0% refactored"

# Row 49
$ws.Range("A49").Value = "50% refactored"

# Row 50
$ws.Range("A50").Value = "100% erroneous code"

# Row 51
$ws.Range("A51").Value = "12% erroneous code"

# Row 46 (filled in after row 51)
$ws.Range("A46").Value = "No coding standard, 0% refactored, 100% erroneous code (everything wrong)"

# Row 52
$ws.Range("A52").Value = "0% refactored, 0% erroneous code"

# Row 53
$ws.Range("A53").Value = "0% refactored, no coding standard, 0% erroneous code"

# Row 54
$ws.Range("A54").Value = "100% refactored, 6% erroneous code"

# Row 55
$ws.Range("A55").Value = "50% refactored, 3% erroneous code"

# Resize column A to fit the new, wider content (target stored width 70.5703125;
# 69.65 is the input that the engine's column-width quantization maps closest to it)
$ws.Columns.Item(1).ColumnWidth = 69.65

# Scroll / selection state to mirror final view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B44").Select()
